$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of roster rows (2-16), values taken from the target state
$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Dereck Lively II", "C", "Dallas Mavericks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
